$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row B7: simple one-decimal format
$ws.Range("B7").Value = 5.6
$ws.Range("B7").NumberFormat = "0.0"

# New row B9: two decimals, negatives in red
$ws.Range("B9").Value = -123
$ws.Range("B9").NumberFormat = "0.00;[Red]0.00"

# New row B11: two decimals, negatives in red with leading minus + trailing space padding
$ws.Range("B11").Value = -123
$ws.Range("B11").NumberFormat = "0.00_ ;[Red]\-0.00\ "

# New row B13: thousands separator with two decimals (built-in format)
$ws.Range("B13").Value = 123000.5
$ws.Range("B13").NumberFormat = "#,##0.00"

# Widen column B to fit the new values
$ws.Columns("B").ColumnWidth = 12.86

# Move the active selection to the last populated cell, like the author did
$ws.Range("B13").Select()
